# Daily attendance processing - 2026-01-18 17:33:26
#
# In the "Recorded By" column (G) of the session analysis sheet, some
# sessions were recorded by both the user and the automated "System"
# recorder, listed as "dnasr281@gmail.com, System". Normalize the order
# of these two names to "System, dnasr281@gmail.com" wherever they occur.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
